# Updates the "cryptos" price/volume table to refresh values as of the
# latest scrape (Tue May 23 06:42:43 UTC 2023). Most rows just get new
# Price (column D) / Volume(1h) (column E) figures, but a handful of rows
# had their rank order swapped (13/14, 32/33, 37/38) -- for those we also
# rewrite the Coin name (B) and Link (C) columns.
#
# Numeric-looking Price strings (e.g. "1.006") are written with a leading
# apostrophe so Excel keeps them as literal text (matching the source
# data's inline-string storage) instead of silently re-parsing them as
# numbers and normalizing away meaningful trailing/precision digits
# (e.g. "18.50" -> 18.5). Values that already fail numeric parsing
# (double-dot separators like "27.424.27", or the " +x.xx% " volume
# strings) are left as plain literals since Excel stores them as text
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "27.424.27"
$ws.Cells.Item(2,5).Value = "  +1.46%  "
$ws.Cells.Item(3,4).Value = "1.864.23"
$ws.Cells.Item(3,5).Value = "  +2.21%  "
$ws.Cells.Item(4,4).Value = "'1.006"
$ws.Cells.Item(4,5).Value = "  -0.23%  "
$ws.Cells.Item(5,4).Value = "'315.26"
$ws.Cells.Item(5,5).Value = "  +1.95%  "
$ws.Cells.Item(6,4).Value = "'1.005"
$ws.Cells.Item(6,5).Value = "  -0.19%  "
$ws.Cells.Item(7,4).Value = "'0.4647"
$ws.Cells.Item(7,5).Value = "  -0.37%  "
$ws.Cells.Item(8,4).Value = "'0.3719"
$ws.Cells.Item(8,5).Value = "  +1.63%  "
$ws.Cells.Item(9,4).Value = "'0.07365"
$ws.Cells.Item(9,5).Value = "  +1.80%  "
$ws.Cells.Item(10,4).Value = "'0.8876"
$ws.Cells.Item(10,5).Value = "  +3.14%  "
$ws.Cells.Item(11,5).Value = "  +4.88%  "
$ws.Cells.Item(12,4).Value = "'19.92"
$ws.Cells.Item(12,5).Value = "  +0.23%  "
$ws.Cells.Item(13,2).Value = "WrappedEther"
$ws.Cells.Item(13,3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13,4).Value = "1.836.81"
$ws.Cells.Item(13,5).Value = "  +3.16%  "
$ws.Cells.Item(14,2).Value = "Polkadot"
$ws.Cells.Item(14,3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(14,4).Value = "'5.409"
$ws.Cells.Item(14,5).Value = "  +1.38%  "
$ws.Cells.Item(15,4).Value = "'6.592"
$ws.Cells.Item(15,5).Value = "  +1.78%  "
$ws.Cells.Item(16,4).Value = "'92.49"
$ws.Cells.Item(16,5).Value = "  +0.70%  "
$ws.Cells.Item(17,5).Value = "  -0.18%  "
$ws.Cells.Item(18,4).Value = "'0.000008896"
$ws.Cells.Item(18,5).Value = "  +2.97%  "
$ws.Cells.Item(19,5).Value = "  -0.09%  "
$ws.Cells.Item(20,4).Value = "'14.88"
$ws.Cells.Item(20,5).Value = "  +2.70%  "
$ws.Cells.Item(21,4).Value = "27.457.60"
$ws.Cells.Item(21,5).Value = "  +2.81%  "
$ws.Cells.Item(22,4).Value = "'5.149"
$ws.Cells.Item(22,5).Value = "  -0.03%  "
$ws.Cells.Item(23,4).Value = "'10.57"
$ws.Cells.Item(23,5).Value = "  +0.35%  "
$ws.Cells.Item(24,4).Value = "2.046.38"
$ws.Cells.Item(24,5).Value = "  +6.43%  "
$ws.Cells.Item(25,4).Value = "'1.904"
$ws.Cells.Item(25,5).Value = "  +3.37%  "
$ws.Cells.Item(26,4).Value = "'153.45"
$ws.Cells.Item(26,5).Value = "  +1.25%  "
$ws.Cells.Item(27,4).Value = "'18.50"
$ws.Cells.Item(27,5).Value = "  +1.73%  "
$ws.Cells.Item(28,4).Value = "'2.077"
$ws.Cells.Item(28,5).Value = "  +0.93%  "
$ws.Cells.Item(29,4).Value = "'5.139"
$ws.Cells.Item(29,5).Value = "  +0.27%  "
$ws.Cells.Item(30,4).Value = "'116.74"
$ws.Cells.Item(30,5).Value = "  +1.30%  "
$ws.Cells.Item(31,4).Value = "'0.08897"
$ws.Cells.Item(31,5).Value = "  +0.65%  "
$ws.Cells.Item(32,2).Value = "HuobiToken"
$ws.Cells.Item(32,3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(32,4).Value = "'3.028"
$ws.Cells.Item(32,5).Value = "  +2.66%  "
$ws.Cells.Item(33,2).Value = "ImmutableX"
$ws.Cells.Item(33,3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(33,4).Value = "'0.7549"
$ws.Cells.Item(33,5).Value = "  +4.84%  "
$ws.Cells.Item(34,4).Value = "'1.164"
$ws.Cells.Item(34,5).Value = "  +2.61%  "
$ws.Cells.Item(35,4).Value = "'4.487"
$ws.Cells.Item(35,5).Value = "  +1.40%  "
$ws.Cells.Item(36,4).Value = "'2.634"
$ws.Cells.Item(36,5).Value = "  +9.44%  "
$ws.Cells.Item(37,2).Value = "TrustWalletToken"
$ws.Cells.Item(37,3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(37,4).Value = "'1.082"
$ws.Cells.Item(37,5).Value = "  +0.06%  "
$ws.Cells.Item(38,2).Value = "VeChain"
$ws.Cells.Item(38,3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(38,4).Value = "'0.01964"
$ws.Cells.Item(38,5).Value = "  +2.09%  "
$ws.Cells.Item(39,4).Value = "'0.05264"
$ws.Cells.Item(39,5).Value = "  +0.07%  "
$ws.Cells.Item(40,4).Value = "'2.988"
$ws.Cells.Item(40,5).Value = "  +2.03%  "
$ws.Cells.Item(41,4).Value = "'7.146"
$ws.Cells.Item(41,5).Value = "  -0.27%  "
$ws.Cells.Item(42,4).Value = "'0.5184"
$ws.Cells.Item(42,5).Value = "  +0.30%  "
$ws.Cells.Item(43,4).Value = "'0.1645"
$ws.Cells.Item(43,5).Value = "  +0.89%  "
$ws.Cells.Item(44,4).Value = "'8.342"
$ws.Cells.Item(44,5).Value = "  +2.00%  "
$ws.Cells.Item(45,4).Value = "'0.4856"
$ws.Cells.Item(45,5).Value = "  +0.90%  "
$ws.Cells.Item(46,4).Value = "'10.39"
$ws.Cells.Item(46,5).Value = "  +2.62%  "
$ws.Cells.Item(47,4).Value = "'1.005"
$ws.Cells.Item(47,5).Value = "  -0.19%  "
$ws.Cells.Item(48,4).Value = "'103.79"
$ws.Cells.Item(48,5).Value = "  +1.11%  "
$ws.Cells.Item(49,4).Value = "'1.652"
$ws.Cells.Item(49,5).Value = "  +1.80%  "
$ws.Cells.Item(50,4).Value = "'0.06244"
$ws.Cells.Item(50,5).Value = "  -0.07%  "
$ws.Cells.Item(51,4).Value = "'65.56"
$ws.Cells.Item(51,5).Value = "  +2.75%  "
